$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-30 down to 21-31.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new data record.
$ws.Cells.Item(20,1).Value = 3
$ws.Cells.Item(20,2).Value = "Femacal de La Calera"
$ws.Cells.Item(20,3).Value = "Coquimbo"
$ws.Cells.Item(20,4).Value = 44839
$ws.Cells.Item(20,5).Value = 5
$ws.Cells.Item(20,6).Value = 300000000
$ws.Cells.Item(20,7).Value = "Espárragos"
$ws.Cells.Item(20,8).Value = "Verde"
$ws.Cells.Item(20,9).Value = "Primera"
$ws.Cells.Item(20,10).Value = 3400
$ws.Cells.Item(20,11).Value = 1400
$ws.Cells.Item(20,12).Value = 1500
$ws.Cells.Item(20,13).Value = 1447
$ws.Cells.Item(20,14).Value = "`$/kilo"
$ws.Cells.Item(20,15).Value = "Provincia de Quillota"
$ws.Cells.Item(20,16).Value = 1447
$ws.Cells.Item(20,17).Value = 1
$ws.Cells.Item(20,18).Value = "Hortaliza"
